$wb = $excel.ActiveWorkbook

# --- pv_costs sheet: add a "Unit" header row under the size-limit row ---
$pv = $wb.Worksheets.Item("pv_costs")

# Shift existing data down by inserting a new row 2
$pv.Rows.Item(2).Insert()

# Populate the new units row
$pv.Range("A2").Value = "Unit"
$pv.Range("B2").Value = "`$/W"
$pv.Range("C2").Value = "`$/W"
$pv.Range("D2").Value = "`$/W"

# Match formatting used by the row above (row 1): centered, wrapped text, row height 16
$pv.Range("A2:D2").HorizontalAlignment = -4108   # xlCenter
$pv.Range("A2:E2").WrapText = $true
$pv.Rows.Item(2).RowHeight = 16

# Make pv_costs the active sheet / selected cell E2 (below the new unit cell)
$pv.Activate()
$pv.Range("E2").Select()
